$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the old "Название доклада" (D) to make room
# for the new "Соавторы" column. This shifts D:G -> E:H.
$ws.Columns.Item(4).Insert()
$ws.Range("D1").Value = "Соавторы"

# Row 2: "Атомарные вычисления" report by Хмельный Никита, with coauthors,
# currently "На модерации", plus the generated report filename.
$ws.Range("A2").Value = "РиМ-2021"
$ws.Range("B2").Value = "Атомарные вычисления"
$ws.Range("C2").Value = "Хмельный Никита "
$ws.Range("D2").Value = "Штопор Александр. Евгеньевич, Торцев Петр. Константинович"
$ws.Range("E2").Value = "ИИ в ИИ"
$ws.Range("H2").Value = "На модерации"
$ws.Range("I2").Value = "Хмельный_Никита_ИИ_в_ИИ_report.doc"

# Row 3: second report by the same author, no coauthors.
$ws.Range("A3").Value = "РиМ-2021"
$ws.Range("B3").Value = "Атомарные вычисления"
$ws.Range("C3").Value = "Хмельный Никита "
$ws.Range("E3").Value = "Вождение пьяным за рулем"
$ws.Range("H3").Value = "На модерации"
$ws.Range("I3").Value = "Хмельный_Никита_Вождение_пьяным_за_рулем_report.doc"
